# Actualización automática 2025-08-22 13:35:10
# Applies the recorded data refresh for "RIOS CARRION ANGEL BENIGNO":
# a new sale of 366.34 for client "CULMA OVIEDO NINI JOHANA" in the
# "240X80 PORCELANATO" group during "agosto", propagated to the three
# report sheets (VENTAS POR GRUPO, VENTA MENSUAL, CUMPLIMIENTO MENSUAL).

$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO -------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# CULMA OVIEDO NINI JOHANA / 240X80 PORCELANATO
$wsGrupo.Range("D10").Value = 366.34

# Summary row: count of advisors with sales in the 240X80 PORCELANATO column
$wsGrupo.Range("D24").Value = "1 de 22"

# --- Sheet: VENTA MENSUAL ----------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# CULMA OVIEDO NINI JOHANA / agosto
$wsMensual.Range("F10").Value = 366.34

# Totals row for agosto
$wsMensual.Range("F24").Value = 3504.01

# --- Sheet: CUMPLIMIENTO MENSUAL ---------------------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# 240X80 PORCELANATO group totals
$wsCumplimiento.Range("D3").Value = 284.93
$wsCumplimiento.Range("E3").Value = 8383.98
$wsCumplimiento.Range("F3").Value = 0.03286803069820773

# Overall TOTAL row
$wsCumplimiento.Range("D19").Value = 3504.01
$wsCumplimiento.Range("E19").Value = 51519.15386304603
$wsCumplimiento.Range("F19").Value = 0.06368245215272543

# Column F narrowed slightly (label column now holds shorter numbers)
$wsCumplimiento.Columns.Item(6).ColumnWidth = 24.17
